$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so cells can be edited, then
# restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A16).
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.03103547731048665
$ws.Range("E2").Value = 0.007414897202561388

$ws.Range("D3").Value = 0.0231563448177947
$ws.Range("E3").Value = 0.005009541984732913

$ws.Range("D4").Value = 0.05164481694847366
$ws.Range("E4").Value = 0.001624883936861776

$ws.Range("D5").Value = 0.1369004412568062
$ws.Range("E5").Value = 0.004056466006814974

$ws.Range("D6").Value = 0.03109917888436214
$ws.Range("E6").Value = 0.001413427561837377

$ws.Range("D7").Value = 0.1179184300171219
$ws.Range("E7").Value = -0.003164556962025333

$ws.Range("D8").Value = 0.1027468776826355
$ws.Range("E8").Value = 0.0007229351165731757

$ws.Range("D9").Value = 0.02935290954000412
$ws.Range("E9").Value = -0.001021033285685036

$ws.Range("D10").Value = 0.1271445804769485
$ws.Range("E10").Value = 0.001109262340543538

$ws.Range("D11").Value = 0.2449243244107252
$ws.Range("E11").Value = 0.002370532458059627

$ws.Range("D12").Value = 0.1040766186546416
$ws.Range("E12").Value = -0.002845759817871496

$ws.Range("E13").Value = 0.0011259429478363

$ws.Protect()
